$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$NL = [char]10

# --- Row 2: Expected Results text got more specific; Actual Results now filled in ---
$ws.Range("F2").Value = "Login success, student home page displayed with correct bidding summary"
$ws.Range("G2").Value = "Login successful but redirected to 'plan bid' instead of 'home' page"

# --- Rows 3-6: Actual Results column filled in (rest of row content unchanged) ---
$ws.Range("G3").Value = "Matched expected results"
$ws.Range("G4").Value = "Matched expected results"
$ws.Range("G5").Value = "Matched expected results"
$ws.Range("G6").Value = "Matched expected results"

# --- New rows 7 and 8: copy formatting (border + wrap text) from the last existing row ---
$ws.Range("A6:G6").Copy()
$ws.Range("A7:G7").PasteSpecial(-4122)
$ws.Range("A6:G6").Copy()
$ws.Range("A8:G8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 7 content
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 1.1
$ws.Range("C7").Value = "Validate that student with correct credentials can login (double check)"
$ws.Range("D7").Value = "Username: ben.ng.2009" + $NL + "Password: qwerty129"
$ws.Range("E7").Value = "Put username into username input field, put password into password input field. Submit"
$ws.Range("F7").Value = "Login success, student home page displayed with correct bidding summary"
$ws.Range("G7").Value = "Login successful but wrong name was displayed on home page"

# Row 8 content
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 1.1
$ws.Range("C8").Value = "Validate that student with correct credentials can login (double double check)"
$ws.Range("D8").Value = "Username: calvin.ng.2009" + $NL + "Password: qwerty130"
$ws.Range("E8").Value = "Put username into username input field, put password into password input field. Submit"
$ws.Range("F8").Value = "Login success, student home page displayed with correct bidding summary"
$ws.Range("G8").Value = "Login successful, bidding summary matches expected result"

# Row heights for the new rows (matching the other short rows in the table)
$ws.Rows.Item(7).RowHeight = 29
$ws.Rows.Item(8).RowHeight = 29

# Matches the active selection saved in the authored workbook
$ws.Range("C4").Select()
